# edit.ps1
# Applies the sem5_timetable.xlsx changes described by the commit:
#   - Section_A / Section_B: collapse the 12-row grid into an 8-row grid with
#     new (re-sequenced) time slots and a strict LTPSC-compliant session layout.
#   - Elective_Coordination: add a "Session Type" column, rename "Slot Name"
#     to "Duration", and expand from a single CS461 row to the full
#     Lecture 1 / Lecture 2 / Tutorial schedule.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Section_A
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

# Remove the now-unused rows 9-12 (grid shrinks from A1:F12 to A1:F8)
$wsA.Rows.Item(9).Delete()
$wsA.Rows.Item(9).Delete()
$wsA.Rows.Item(9).Delete()
$wsA.Rows.Item(9).Delete()

$wsA.Range("A2").Value = "09:00-10:30"
$wsA.Range("B2").Value = "Free"
$wsA.Range("C2").Value = "CS303"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "CS303"
$wsA.Range("F2").Value = "CS304"

$wsA.Range("A3").Value = "10:30-12:00"
$wsA.Range("B3").Value = "CS304"
$wsA.Range("C3").Value = "CS309"
$wsA.Range("D3").Value = "CS309"
$wsA.Range("E3").Value = "Free"
$wsA.Range("F3").Value = "CS461 (Elective)"

$wsA.Range("A4").Value = "12:00-13:00"
$wsA.Range("B4").Value = "LUNCH BREAK"
$wsA.Range("C4").Value = "LUNCH BREAK"
$wsA.Range("D4").Value = "LUNCH BREAK"
$wsA.Range("E4").Value = "LUNCH BREAK"
$wsA.Range("F4").Value = "LUNCH BREAK"

$wsA.Range("A5").Value = "13:00-14:30"
$wsA.Range("B5").Value = "CS309"
$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "Free"
$wsA.Range("E5").Value = "CS461 (Elective)"
$wsA.Range("F5").Value = "CS303"

$wsA.Range("A6").Value = "14:30-15:30"
$wsA.Range("B6").Value = "CS303 (Tutorial)"
$wsA.Range("C6").Value = "CS304 (Tutorial)"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "Free"
$wsA.Range("F6").Value = "CS461 (Tutorial)"

$wsA.Range("A7").Value = "15:30-17:00"
$wsA.Range("B7").Value = "Free"
$wsA.Range("C7").Value = "Free"
$wsA.Range("D7").Value = "Free"
$wsA.Range("E7").Value = "CS304"
$wsA.Range("F7").Value = "Free"

$wsA.Range("A8").Value = "17:00-18:00"
$wsA.Range("B8").Value = "Free"
$wsA.Range("C8").Value = "Free"
$wsA.Range("D8").Value = "Free"
$wsA.Range("E8").Value = "CS309 (Tutorial)"
$wsA.Range("F8").Value = "Free"

# ---------------------------------------------------------------------------
# Section_B
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Rows.Item(9).Delete()
$wsB.Rows.Item(9).Delete()
$wsB.Rows.Item(9).Delete()
$wsB.Rows.Item(9).Delete()

$wsB.Range("A2").Value = "09:00-10:30"
$wsB.Range("B2").Value = "Free"
$wsB.Range("C2").Value = "CS309"
$wsB.Range("D2").Value = "CS304"
$wsB.Range("E2").Value = "CS303"
$wsB.Range("F2").Value = "CS304"

$wsB.Range("A3").Value = "10:30-12:00"
$wsB.Range("B3").Value = "CS309"
$wsB.Range("C3").Value = "Free"
$wsB.Range("D3").Value = "Free"
$wsB.Range("E3").Value = "Free"
$wsB.Range("F3").Value = "CS461 (Elective)"

$wsB.Range("A4").Value = "12:00-13:00"
$wsB.Range("B4").Value = "LUNCH BREAK"
$wsB.Range("C4").Value = "LUNCH BREAK"
$wsB.Range("D4").Value = "LUNCH BREAK"
$wsB.Range("E4").Value = "LUNCH BREAK"
$wsB.Range("F4").Value = "LUNCH BREAK"

$wsB.Range("A5").Value = "13:00-14:30"
$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "CS303"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "CS461 (Elective)"
$wsB.Range("F5").Value = "Free"

$wsB.Range("A6").Value = "14:30-15:30"
$wsB.Range("B6").Value = "CS309 (Tutorial)"
$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "Free"
$wsB.Range("E6").Value = "Free"
$wsB.Range("F6").Value = "CS461 (Tutorial)"

$wsB.Range("A7").Value = "15:30-17:00"
$wsB.Range("B7").Value = "CS303"
$wsB.Range("C7").Value = "CS304"
$wsB.Range("D7").Value = "Free"
$wsB.Range("E7").Value = "Free"
$wsB.Range("F7").Value = "CS309"

$wsB.Range("A8").Value = "17:00-18:00"
$wsB.Range("B8").Value = "Free"
$wsB.Range("C8").Value = "Free"
$wsB.Range("D8").Value = "Free"
$wsB.Range("E8").Value = "CS304 (Tutorial)"
$wsB.Range("F8").Value = "Free"

# ---------------------------------------------------------------------------
# Elective_Coordination
# ---------------------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Elective_Coordination")

# Insert a new "Session Type" column before the existing "Day" column (B).
# This shifts Day -> C, Time Slot -> D, Slot Name -> E, Sections -> F.
$wsE.Columns.Item(2).Insert()

$wsE.Range("B1").Value = "Session Type"
$wsE.Range("E1").Value = "Duration"

$wsE.Range("A2").Value = "CS461"
$wsE.Range("B2").Value = "Lecture 1"
$wsE.Range("C2").Value = "Fri"
$wsE.Range("D2").Value = "10:30-12:00"
$wsE.Range("E2").Value = "1.5 hours"
$wsE.Range("F2").Value = "A & B (Common Slot)"

$wsE.Range("A3").Value = "CS461"
$wsE.Range("B3").Value = "Lecture 2"
$wsE.Range("C3").Value = "Thu"
$wsE.Range("D3").Value = "13:00-14:30"
$wsE.Range("E3").Value = "1.5 hours"
$wsE.Range("F3").Value = "A & B (Common Slot)"

$wsE.Range("A4").Value = "CS461"
$wsE.Range("B4").Value = "Tutorial"
$wsE.Range("C4").Value = "Fri"
$wsE.Range("D4").Value = "14:30-15:30"
$wsE.Range("E4").Value = "1 hour"
$wsE.Range("F4").Value = "A & B (Common Slot)"
